$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every touched cell so numeric-looking values
# (e.g. "1.000", "125.00") are preserved verbatim, matching the sheet's
# existing all-text data model, instead of Excel auto-coercing them to numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.038.08'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.72%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.766.69'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.73%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.71'

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9993'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.56%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4271'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.45%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3623'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -4.63%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.64'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.80%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07455'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.67%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.096'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.37%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.75%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.15'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.42%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.076'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -4.36%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.328'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.11%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.773.29'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.13%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.34'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.46%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001059'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.94%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06391'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.08%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.70%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.08'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.13%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.976'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.64%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.030.51'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.85%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.27%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.130'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.75%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.30'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.97%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.20'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.94%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.980.14'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.97%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.143'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -10.66%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.00'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.98%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.161'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.63%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.629'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.46%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.645'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.84%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08884'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.53%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.28%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02310'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.66%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2109'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.52%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.018'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.00%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06050'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.86%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6361'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.75%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.184'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.03%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9990'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.55%  '

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'WEMIXTOKEN'

$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.403'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.05%  '

$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'FraxShare'

$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.834'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.76%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.48'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.37%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5923'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.62%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.680'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.32%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.005'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.42%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '122.21'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.32%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.180'

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.10%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06864'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.34%  '
